$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2" "303.91"
Set-TextValue "E2" "5.21%"
Set-TextValue "D3" "35.31"
Set-TextValue "E3" "13.98%"
Set-TextValue "D4" "5.195"
Set-TextValue "E4" "4.89%"
Set-TextValue "D5" "0.07821"
Set-TextValue "E5" "6.04%"
Set-TextValue "D6" "2.289"
Set-TextValue "E6" "-1.54%"
Set-TextValue "D7" "8.030"
Set-TextValue "E7" "4.03%"
Set-TextValue "D8" "3.984"
Set-TextValue "E8" "7.13%"
Set-TextValue "D9" "0.9286"
Set-TextValue "E9" "1.96%"
Set-TextValue "E10" "10.37%"
Set-TextValue "D11" "0.1840"
Set-TextValue "E11" "8.49%"
Set-TextValue "D12" "0.08581"
Set-TextValue "E12" "3.77%"
Set-TextValue "D13" "0.03383"
Set-TextValue "E13" "8.42%"
Set-TextValue "D14" "0.09916"
Set-TextValue "E14" "-0.46%"
Set-TextValue "D15" "0.001479"
Set-TextValue "E15" "-0.87%"
Set-TextValue "D16" "0.04650"
Set-TextValue "E16" "3.31%"
Set-TextValue "D17" "0.005754"
Set-TextValue "E17" "-0.69%"
Set-TextValue "D18" "3.465"
Set-TextValue "E18" "-0.97%"
Set-TextValue "D19" "2.111"
Set-TextValue "E19" "-0.04%"
Set-TextValue "D20" "0.3419"
Set-TextValue "E20" "2.83%"
Set-TextValue "D21" "0.1325"
Set-TextValue "E21" "3.07%"
Set-TextValue "D22" "4.560"
Set-TextValue "E22" "9.02%"
Set-TextValue "E23" "13.62%"
Set-TextValue "D24" "0.001223"
Set-TextValue "E24" "1.25%"
Set-TextValue "D25" "0.004447"
Set-TextValue "E25" "6.54%"
Set-TextValue "D26" "0.0001299"
Set-TextValue "E26" "0.15%"
Set-TextValue "D27" "0.0003397"
Set-TextValue "E27" "0.25%"
Set-TextValue "D39" "0.01750"
Set-TextValue "E39" "10.70%"
Set-TextValue "D40" "0.04732"
Set-TextValue "E40" "5.99%"
Set-TextValue "D41" "0.007693"
Set-TextValue "E41" "4.72%"
Set-TextValue "D42" "0.1412"
Set-TextValue "E42" "6.37%"
Set-TextValue "D43" "0.007046"
Set-TextValue "E43" "-25.73%"
Set-TextValue "D44" "0.002299"
Set-TextValue "E44" "1.93%"
Set-TextValue "D45" "0.009924"
Set-TextValue "E45" "22.99%"
Set-TextValue "D46" "0.00005999"
Set-TextValue "E46" "-1.67%"
Set-TextValue "D47" "0.00000000749"
Set-TextValue "E47" "0.08%"
Set-TextValue "D48" "5.812"
Set-TextValue "E48" "126.58%"
Set-TextValue "D49" "0.002688"
Set-TextValue "E49" "34.61%"
Set-TextValue "D50" "0.00002098"
Set-TextValue "E50" "0.08%"
Set-TextValue "D51" "0.0001998"
Set-TextValue "E51" "0.08%"
